# Update "想去人数" (want-to-go count) figures in the 展览 sheet and the
# mirrored rows in the 全部类型 sheet (which aggregates all event types).
#
# 展览 is worksheet #1, 全部类型 is worksheet #4 (see workbook.xml sheet order:
# 展览, 演出, 本地生活, 全部类型).

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item(1)   # 展览
$wsAll     = $wb.Worksheets.Item(4)   # 全部类型

# Map of column-F cell -> new value for the 展览 sheet.
$exhibitUpdates = @{
    "F3"  = 1233
    "F4"  = 16969
    "F5"  = 38
    "F6"  = 1659
    "F7"  = 72
    "F9"  = 11
    "F13" = 11744
    "F15" = 14
    "F16" = 1433
    "F17" = 4670
    "F18" = 479
    "F20" = 409
    "F22" = 905
    "F25" = 31
}

foreach ($cell in $exhibitUpdates.Keys) {
    $wsExhibit.Range($cell).Value = $exhibitUpdates[$cell]
}

# Same events, mirrored at different rows inside the aggregate 全部类型 sheet.
$allUpdates = @{
    "F4"  = 1233
    "F5"  = 16969
    "F6"  = 38
    "F7"  = 1659
    "F8"  = 72
    "F10" = 11
    "F16" = 11744
    "F18" = 14
    "F19" = 1433
    "F20" = 4670
    "F21" = 479
    "F23" = 409
    "F25" = 905
    "F28" = 31
}

foreach ($cell in $allUpdates.Keys) {
    $wsAll.Range($cell).Value = $allUpdates[$cell]
}
